$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "57×41=2337" "68×82=5576"
Replace-Text "94×74=6956" "50×74=3700"
Replace-Text "81×54=4374" "64×42=2688"
Replace-Text "61×87=5307" "34×51=1734"
Replace-Text "79×41=3239" "51×91=4641"
Replace-Text "66×98=6468" "76×27=2052"
Replace-Text "40×40=1600" "83×72=5976"
Replace-Text "23×86=1978" "41×14=574"
Replace-Text "97×83=8051" "17×35=595"
Replace-Text "74×82=6068" "16×46=736"
Replace-Text "31×37=1147" "21×41=861"
Replace-Text "27×76=2052" "25×76=1900"
Replace-Text "61×75=4575" "87×69=6003"
Replace-Text "69×82=5658" "71×60=4260"
Replace-Text "85×86=7310" "81×75=6075"
Replace-Text "58×91=5278" "72×85=6120"
Replace-Text "45×17=765" "36×50=1800"
Replace-Text "57×22=1254" "26×32=832"
Replace-Text "28×89=2492" "49×67=3283"
Replace-Text "85×53=4505" "39×67=2613"
Replace-Text "24×70=1680" "53×69=3657"
Replace-Text "96×93=8928" "13×18=234"
Replace-Text "67×44=2948" "55×52=2860"
Replace-Text "39×61=2379" "55×96=5280"
Replace-Text "20×12=240" "56×89=4984"
